# v03-slides-and-speaker-notes.pptx -- "More work on 5505 video03"
#
# Slide 27 ("White-space delimited files"): switch the read_delim() call to
# read_table() and drop the now-unused `delim=NULL,` argument (and the line
# break that introduced it), then update the warning message that the new
# call produces.
#
# Slide 28 ("Display the white-space data"): update the printed tibble to
# reflect the new (correct) 2-column parse result.
#
# NOTE: this runtime's TextRange/Characters API is unreliable when an edit's
# start/end lands exactly on a soft line-break (chr(11), inserted for
# <a:br/>) -- the break silently survives the edit. Keeping at least one
# ordinary character on *both* sides of any break that's inside a replaced
# range avoids that, so replacements below are deliberately padded with a
# bit of unchanged context text and Characters()/IndexOf() (not Run.Text,
# whose Start/Length are unreliable once a paragraph has >1 run) is used
# throughout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 27: fn <- ...; raw_data <- read_delim(fn, delim=NULL, col_names=...)
# ---------------------------------------------------------------------
$slide27 = $p.Slides.Item(27)
$code27 = $slide27.Shapes.Item(2).TextFrame.TextRange

# 1) read_delim -> read_table
$full = $code27.Text
$old = "read_delim"
$idx = $full.IndexOf($old)
$code27.Characters($idx + 1, $old.Length).Text = "read_table"

# 2) drop the "  delim=NULL," argument together with the line break that
#    precedes it, e.g.
#      read_table(fn, <br/>  delim=NULL,<br/>  col_names=...
#    becomes
#      read_table(fn, <br/>  col_names=...
#    The replaced span starts one character early (the trailing space of
#    "(fn, ") so the leading <br/> isn't the very first character of the
#    edit (see note above); that character is simply written back unchanged.
$full = $code27.Text
$old = "(fn, `v  delim=NULL,"
$idx = $full.IndexOf($old)
$code27.Characters($idx + 1, $old.Length).Text = "(fn, "

# 3) Updated warning text (this paragraph is a single run, so replacing the
#    whole thing is safe).
$full = $code27.Text
$old = "## Warning: One or more parsing issues, see`n## ``problems()`` for details"
$idx = $full.IndexOf($old)
$code27.Characters($idx + 1, $old.Length).Text = "## Warning: Insufficient ``col_types``. Guessing 1`n## columns."

# ---------------------------------------------------------------------
# Slide 28: raw_data; ## # A tibble: 4 x 4 ...
# ---------------------------------------------------------------------
$slide28 = $p.Slides.Item(28)
$code28 = $slide28.Shapes.Item(2).TextFrame.TextRange

$full = $code28.Text
$old = "## # A tibble: 4 x 4`n" + `
       "##       x y     X3        X4`n" + `
       "##   <dbl> <lgl> <chr>  <dbl>`n" + `
       "## 1     1 NA     <NA>      4`n" + `
       "## 2     2 NA     <NA>      8`n" + `
       "## 3     3 NA    ""12\r""    NA`n" + `
       "## 4     4 NA    ""16\r""    NA"
$idx = $full.IndexOf($old)
$new = "## # A tibble: 4 x 2`n" + `
       "##       x     y`n" + `
       "##   <dbl> <dbl>`n" + `
       "## 1     1     4`n" + `
       "## 2     2     8`n" + `
       "## 3     3    12`n" + `
       "## 4     4    16"
$code28.Characters($idx + 1, $old.Length).Text = $new
